# Insert a new data row before row 316 (weekly Fruta/Hortaliza update),
# shifting the existing rows 316-366 down to 317-367, then populate the
# newly-inserted row 316 with the latest "Polar King" Durazno record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 316..366 down to 317..367.
$ws.Rows.Item(316).Insert()

# Fill in the new row 316 with the new record.
$ws.Cells.Item(316, 1).Value  = 11
$ws.Cells.Item(316, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(316, 3).Value  = "Bíobío"
$ws.Cells.Item(316, 4).Value  = 45258
$ws.Cells.Item(316, 5).Value  = 8
$ws.Cells.Item(316, 6).Value  = "Fruta"
$ws.Cells.Item(316, 7).Value  = 100103
$ws.Cells.Item(316, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(316, 9).Value  = 100103004
$ws.Cells.Item(316, 10).Value = "Durazno"
$ws.Cells.Item(316, 11).Value = "Polar King"
$ws.Cells.Item(316, 12).Value = "Primera"
$ws.Cells.Item(316, 13).Value = 100
$ws.Cells.Item(316, 14).Value = 14000
$ws.Cells.Item(316, 15).Value = 15000
$ws.Cells.Item(316, 16).Value = 14500
$ws.Cells.Item(316, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(316, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(316, 19).Value = 967
$ws.Cells.Item(316, 20).Value = 15
